# Temporary Fix for RestDayEquality Penalty Value
# - "Quality" column header renamed to "Penalty"
# - Recomputed penalty-related figures for both the Pre- and
#   Post-Optimization tables (games counts, timeslot/game ratios,
#   success percentages and the weighted penalty totals) now that
#   exhibition games are excluded from the penalty calculation.
# - Optimization Time updated to reflect the faster run.
# - Cells whose inputs feed the RestDayEquality penalty are highlighted
#   (orange shades in the Pre-Optimization table, green shades in the
#   Post-Optimization table; darker shade = primary penalty inputs,
#   lighter shade = secondary contributing inputs).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Colors used to flag the cells feeding the penalty calculation ---
$darkOrange  = 7199231    # RGB(FFD96D) - primary pre-optimization inputs
$lightOrange = 11528959   # RGB(FFEAAF) - secondary pre-optimization inputs
$darkGreen   = 8440478    # RGB(9ECA80) - primary post-optimization inputs
$lightGreen  = 11854022   # RGB(C6E0B4) - secondary post-optimization inputs

# --- Header rename: "Quality" -> "Penalty" (shared by both tables) ---
$ws.Range("K1").Value = "Penalty"
$ws.Range("K15").Value = "Penalty"

# =========================================================
# Schedule (Pre-Optimization) table - rows 2-4
# =========================================================

# Division U7 (Tier: 1) - row 2
$ws.Range("C2").Value = 164
$ws.Range("D2").Value = 63
$ws.Range("H2").Value = "'72.25%"
$ws.Range("J2").Value = 10
$ws.Range("K2").Value = 509.45

# Division U8 (Tier: 0) - row 3
$ws.Range("C3").Value = 160
$ws.Range("D3").Value = 100
$ws.Range("F3").Value = 16
$ws.Range("G3").Value = "'90.91%"
$ws.Range("H3").Value = "'61.54%"
$ws.Range("J3").Value = 10
$ws.Range("K3").Value = 709.27

# Division U9 (Tier: 0) - row 4
$ws.Range("C4").Value = 205
$ws.Range("D4").Value = 71
$ws.Range("H4").Value = "'74.28%"
$ws.Range("J4").Value = 8
$ws.Range("K4").Value = 851.87

# =========================================================
# Schedule (Post-Optimization) table - rows 16-18
# =========================================================

# Division U7 (Tier: 1) - row 16
$ws.Range("C16").Value = 203
$ws.Range("D16").Value = 24
$ws.Range("H16").Value = "'89.43%"
$ws.Range("J16").Value = 12
$ws.Range("K16").Value = 307.45

# Division U8 (Tier: 0) - row 17
$ws.Range("C17").Value = 205
$ws.Range("D17").Value = 55
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = "'91.48%"
$ws.Range("H17").Value = "'78.85%"
$ws.Range("J17").Value = 12
$ws.Range("K17").Value = 454.57

# Division U9 (Tier: 0) - row 18
$ws.Range("C18").Value = 246
$ws.Range("D18").Value = 30
$ws.Range("H18").Value = "'89.13%"
$ws.Range("J18").Value = 10
$ws.Range("K18").Value = 741.87

# --- Optimization Time ---
$ws.Range("A30").Value = "0 min, 11 sec"

# =========================================================
# Highlight the cells that feed the penalty calculation
# =========================================================

# Pre-Optimization - primary inputs (dark orange)
foreach ($addr in @("C2","D2","H2","K2","C3","D3","H3","K3","C4","D4","K4")) {
    $ws.Range($addr).Interior.Color = $darkOrange
}

# Pre-Optimization - secondary inputs (light orange)
foreach ($addr in @("J2","F3","G3","J3","H4","J4")) {
    $ws.Range($addr).Interior.Color = $lightOrange
}

# Post-Optimization - primary inputs (dark green)
foreach ($addr in @("C16","D16","H16","K16","C17","D17","H17","K17","C18","D18","K18")) {
    $ws.Range($addr).Interior.Color = $darkGreen
}

# Post-Optimization - secondary inputs (light green)
foreach ($addr in @("J16","F17","G17","J17","H18","J18")) {
    $ws.Range($addr).Interior.Color = $lightGreen
}
